$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 826.1667  # H5
$ws.Cells.Item(5, 9).Value = 395.83334  # I5
$ws.Cells.Item(5, 10).Value = 1256.5  # J5
$ws.Cells.Item(5, 11).Value = 395.83334  # K5
$ws.Cells.Item(5, 12).Value = 1256.5  # L5
$ws.Cells.Item(5, 13).Value = -280.83334  # M5
$ws.Cells.Item(5, 14).Value = -1486.5  # N5
$ws.Cells.Item(55, 8).Value = 131.66667  # H55
$ws.Cells.Item(62, 8).Value = 14295901  # H62
$ws.Cells.Item(62, 10).Value = 6966.3335  # J62
$ws.Cells.Item(62, 12).Value = 6966.3335  # L62
$ws.Cells.Item(62, 14).Value = -8214.333500000001  # N62
$ws.Cells.Item(65, 8).Value = 14295901  # H65
$ws.Cells.Item(65, 10).Value = 6966.3335  # J65
$ws.Cells.Item(65, 12).Value = 34831.6675  # L65
$ws.Cells.Item(65, 14).Value = -41071.6675  # N65
$ws.Cells.Item(129, 8).Value = 22222846  # H129
$ws.Cells.Item(129, 9).Value = 525.6923  # I129
$ws.Cells.Item(129, 10).Value = 166667920  # J129
$ws.Cells.Item(129, 11).Value = 1577.0769  # K129
$ws.Cells.Item(129, 12).Value = 500003760  # L129
$ws.Cells.Item(129, 13).Value = 3422.9231  # M129
$ws.Cells.Item(129, 14).Value = -500013760  # N129
$ws.Cells.Item(135, 8).Value = 3592.9622  # H135
$ws.Cells.Item(135, 9).Value = 1456.2  # I135
$ws.Cells.Item(135, 11).Value = 13105.8  # K135
$ws.Cells.Item(135, 13).Value = -10570.8  # M135
$ws.Cells.Item(137, 8).Value = 9537.299999999999  # H137
$ws.Cells.Item(137, 9).Value = 5318.65  # I137
$ws.Cells.Item(137, 10).Value = 17974.6  # J137
$ws.Cells.Item(137, 11).Value = 15955.95  # K137
$ws.Cells.Item(137, 12).Value = 53923.8  # L137
$ws.Cells.Item(137, 13).Value = -13405.95  # M137
$ws.Cells.Item(137, 14).Value = -59023.8  # N137
$ws.Cells.Item(138, 8).Value = 3409.1538  # H138
$ws.Cells.Item(138, 9).Value = 1106.32  # I138
$ws.Cells.Item(138, 10).Value = 5541.407  # J138
$ws.Cells.Item(138, 11).Value = 3318.96  # K138
$ws.Cells.Item(138, 12).Value = 16624.221  # L138
$ws.Cells.Item(138, 13).Value = 1821.04  # M138
$ws.Cells.Item(138, 14).Value = -26904.221  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(58, 8).Value = 0  # H58
$ws.Cells.Item(58, 10).Value = 0  # J58
$ws.Cells.Item(58, 12).Value = 0  # L58
$ws.Cells.Item(58, 14).ClearContents()  # N58 removed
$ws.Cells.Item(61, 8).Value = 9823.444  # H61
$ws.Cells.Item(61, 9).Value = 14474.5  # I61
$ws.Cells.Item(61, 11).Value = 14474.5  # K61
$ws.Cells.Item(61, 13).Value = -14262.5  # M61
$ws.Cells.Item(74, 8).Value = 4324.263  # H74
$ws.Cells.Item(74, 9).Value = 4635.3125  # I74
$ws.Cells.Item(74, 11).Value = 4635.3125  # K74
$ws.Cells.Item(74, 13).Value = -3761.3125  # M74
$ws.Cells.Item(77, 8).Value = 4324.263  # H77
$ws.Cells.Item(77, 9).Value = 4635.3125  # I77
$ws.Cells.Item(77, 11).Value = 23176.5625  # K77
$ws.Cells.Item(77, 13).Value = -18808.5625  # M77
$ws.Cells.Item(80, 8).Value = 28394.5  # H80
$ws.Cells.Item(80, 10).Value = 28394.5  # J80
$ws.Cells.Item(80, 12).Value = 28394.5  # L80
$ws.Cells.Item(80, 14).Value = -30390.5  # N80
$ws.Cells.Item(83, 8).Value = 28394.5  # H83
$ws.Cells.Item(83, 10).Value = 28394.5  # J83
$ws.Cells.Item(83, 12).Value = 85183.5  # L83
$ws.Cells.Item(83, 14).Value = -95167.5  # N83
$ws.Cells.Item(136, 8).Value = 9823.444  # H136
$ws.Cells.Item(136, 9).Value = 14474.5  # I136
$ws.Cells.Item(136, 11).Value = 43423.5  # K136
$ws.Cells.Item(136, 13).Value = -40873.5  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 6803383.5  # H20
$ws.Cells.Item(20, 9).Value = 10204622  # I20
$ws.Cells.Item(20, 11).Value = 10204622  # K20
$ws.Cells.Item(20, 13).Value = -10204375  # M20
$ws.Cells.Item(22, 8).Value = 11235.6  # H22
$ws.Cells.Item(22, 9).Value = 8132.185  # I22
$ws.Cells.Item(22, 11).Value = 8132.185  # K22
$ws.Cells.Item(22, 13).Value = -7959.185  # M22
$ws.Cells.Item(134, 8).Value = 3345204  # H134
$ws.Cells.Item(134, 9).Value = 3582182.8  # I134
$ws.Cells.Item(134, 11).Value = 10746548.4  # K134
$ws.Cells.Item(134, 13).Value = -10744013.4  # M134
$ws.Cells.Item(140, 8).Value = 75000  # H140
$ws.Cells.Item(140, 10).Value = 75000  # J140
$ws.Cells.Item(140, 12).Value = 75000  # L140
$ws.Cells.Item(140, 14).Value = -85360  # N140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4683.9165  # H31
$ws.Cells.Item(31, 9).Value = 2310.7144  # I31
$ws.Cells.Item(31, 11).Value = 2310.7144  # K31
$ws.Cells.Item(31, 13).Value = -2015.7144  # M31
$ws.Cells.Item(34, 8).Value = 4683.9165  # H34
$ws.Cells.Item(34, 9).Value = 2310.7144  # I34
$ws.Cells.Item(34, 11).Value = 2310.7144  # K34
$ws.Cells.Item(34, 13).Value = -2108.7144  # M34
$ws.Cells.Item(134, 8).Value = 62509892  # H134
$ws.Cells.Item(134, 9).Value = 90917780  # I134
$ws.Cells.Item(134, 11).Value = 272753340  # K134
$ws.Cells.Item(134, 13).Value = -272750805  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 14287.2  # H56
$ws.Cells.Item(56, 9).Value = 14287.2  # I56
$ws.Cells.Item(56, 11).Value = 14287.2  # K56
$ws.Cells.Item(56, 13).Value = -13757.2  # M56
$ws.Cells.Item(62, 8).Value = 19988.666  # H62
$ws.Cells.Item(62, 10).Value = 19988.666  # J62
$ws.Cells.Item(62, 12).Value = 59965.99800000001  # L62
$ws.Cells.Item(62, 14).Value = -61337.99800000001  # N62
$ws.Cells.Item(65, 8).Value = 19988.666  # H65
$ws.Cells.Item(65, 10).Value = 19988.666  # J65
$ws.Cells.Item(65, 12).Value = 179897.994  # L65
$ws.Cells.Item(65, 14).Value = -186761.994  # N65
$ws.Cells.Item(113, 8).Value = 2501027  # H113
$ws.Cells.Item(113, 9).Value = 6667020.5  # I113
$ws.Cells.Item(113, 10).Value = 1430.8  # J113
$ws.Cells.Item(113, 11).Value = 20001061.5  # K113
$ws.Cells.Item(113, 12).Value = 4292.4  # L113
$ws.Cells.Item(113, 13).Value = -19998891.5  # M113
$ws.Cells.Item(113, 14).Value = -8632.4  # N113
$ws.Cells.Item(131, 8).Value = 39398144  # H131
$ws.Cells.Item(131, 9).Value = 48489940  # I131
$ws.Cells.Item(131, 11).Value = 145469820  # K131
$ws.Cells.Item(131, 13).Value = -145464780  # M131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 63  # H2
$ws.Cells.Item(2, 9).Value = 52  # I2
$ws.Cells.Item(2, 10).Value = 85  # J2
$ws.Cells.Item(2, 11).Value = 52  # K2
$ws.Cells.Item(2, 12).Value = 85  # L2
$ws.Cells.Item(2, 13).Value = 61  # M2
$ws.Cells.Item(2, 14).Value = -311  # N2
$ws.Cells.Item(49, 8).Value = 37500  # H49
$ws.Cells.Item(49, 10).Value = 37500  # J49
$ws.Cells.Item(49, 12).Value = 37500  # L49
$ws.Cells.Item(49, 14).Value = -37868  # N49
$ws.Cells.Item(70, 8).Value = 21785.715  # H70
$ws.Cells.Item(70, 9).Value = 7500  # I70
$ws.Cells.Item(70, 10).Value = 32500  # J70
$ws.Cells.Item(70, 11).Value = 7500  # K70
$ws.Cells.Item(70, 12).Value = 32500  # L70
$ws.Cells.Item(70, 13).Value = -7230  # M70
$ws.Cells.Item(70, 14).Value = -33040  # N70
$ws.Cells.Item(73, 8).Value = 21785.715  # H73
$ws.Cells.Item(73, 9).Value = 7500  # I73
$ws.Cells.Item(73, 10).Value = 32500  # J73
$ws.Cells.Item(73, 11).Value = 7500  # K73
$ws.Cells.Item(73, 12).Value = 32500  # L73
$ws.Cells.Item(73, 13).Value = -6564  # M73
$ws.Cells.Item(73, 14).Value = -34372  # N73
$ws.Cells.Item(109, 8).Value = 40000  # H109
$ws.Cells.Item(109, 10).Value = 40000  # J109
$ws.Cells.Item(109, 12).Value = 40000  # L109
$ws.Cells.Item(109, 14).Value = -42080  # N109
$ws.Cells.Item(112, 8).Value = 74000  # H112
$ws.Cells.Item(112, 10).Value = 74000  # J112
$ws.Cells.Item(112, 12).Value = 74000  # L112
$ws.Cells.Item(112, 14).Value = -76216  # N112
$ws.Cells.Item(132, 8).Value = 25645098  # H132
$ws.Cells.Item(132, 9).Value = 40003944  # I132
$ws.Cells.Item(132, 10).Value = 4306.5  # J132
$ws.Cells.Item(132, 11).Value = 120011832  # K132
$ws.Cells.Item(132, 12).Value = 12919.5  # L132
$ws.Cells.Item(132, 13).Value = -120009302  # M132
$ws.Cells.Item(132, 14).Value = -17979.5  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 892.3333  # H22
$ws.Cells.Item(22, 9).Value = 794.6667  # I22
$ws.Cells.Item(22, 11).Value = 794.6667  # K22
$ws.Cells.Item(22, 13).Value = -499.6667  # M22
$ws.Cells.Item(27, 8).Value = 892.3333  # H27
$ws.Cells.Item(27, 9).Value = 794.6667  # I27
$ws.Cells.Item(27, 11).Value = 794.6667  # K27
$ws.Cells.Item(27, 13).Value = -687.6667  # M27
$ws.Cells.Item(122, 8).Value = 2924.6  # H122
$ws.Cells.Item(122, 9).Value = 3242  # I122
$ws.Cells.Item(122, 10).Value = 2713  # J122
$ws.Cells.Item(122, 11).Value = 9726  # K122
$ws.Cells.Item(122, 12).Value = 8139  # L122
$ws.Cells.Item(122, 13).Value = -7276  # M122
$ws.Cells.Item(122, 14).Value = -13039  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2494.4644  # H122
$ws.Cells.Item(122, 9).Value = 2194.44  # I122
$ws.Cells.Item(122, 11).Value = 6583.32  # K122
$ws.Cells.Item(122, 13).Value = -4133.32  # M122
$ws.Cells.Item(132, 8).Value = 4202.5  # H132
$ws.Cells.Item(132, 9).Value = 4025.9143  # I132
$ws.Cells.Item(132, 11).Value = 12077.7429  # K132
$ws.Cells.Item(132, 13).Value = -9547.742899999999  # M132
$ws.Cells.Item(133, 8).Value = 86502.39999999999  # H133
$ws.Cells.Item(133, 10).Value = 86502.39999999999  # J133
$ws.Cells.Item(133, 12).Value = 86502.39999999999  # L133
$ws.Cells.Item(133, 14).Value = -96622.39999999999  # N133
